# Automatische test-sync: 2025-08-02 00:14:50
#
# Appends a new test-mail log entry (row 8) to the "Logs" sheet and updates
# the corresponding category count (row 6) on the "Dashboard" sheet. The
# bar chart on Dashboard already references the category range, so Excel
# will pick up the new row once the chart's series formulas are extended.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- Logs sheet: append new row 8 -----------------------------------------
$logs.Range("A8").Value = "Bestel je 200 stuks M8-bouten RVS voor Van Dijk?"
$logs.Range("B8").Value = "mailmind.test@zohomail.eu"
$logs.Range("C8").Value = "Testmail #18: Bestel je 200 stuks M8-bouten RVS voor Van Dijk?"
$logs.Range("D8").Value = "Bestelling / Levering"
$logs.Range("E8").Value = "Bedankt, we hebben dit doorgestuurd naar support@bedrijf.nl."
$logs.Range("F8").Value = "2025-08-02 00:14:30"
$logs.Range("G8").Value = "Ja"
$logs.Range("H8").Value = "Ja"
$logs.Range("I8").Value = "Nee"
$logs.Range("J8").Value = "Nee"

# --- extend the conditional-formatting ranges on Logs so row 8 is covered -
# (each block's rules all share one sqref - modifying any single rule's
# applies-to range re-seats the whole block, so one call per column suffices)
$logs.Range("D2:D7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D8"))
$logs.Range("G2:G7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G8"))
$logs.Range("H2:H7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H8"))
$logs.Range("I2:I7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I8"))
$logs.Range("J2:J7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J8"))

# --- Dashboard sheet: append new row 6 for the new category ---------------
$dashboard.Range("A6").Value = "Bestelling / Levering"
$dashboard.Range("B6").Value = 1

# --- Chart: extend the category/value series ranges to include row 6 -----
# Set the raw SERIES() formula (rather than assigning XValues/Values range
# objects) so the stored reference keeps the quoted sheet name + absolute
# addressing the workbook already used ('Dashboard'!$A$2:$A$6 etc.), and the
# series-name reference (B1, no $) is left exactly as it was.
$chart = $dashboard.ChartObjects(1).Chart
$series = $chart.SeriesCollection(1)
$series.Formula = "=SERIES('Dashboard'!B1,'Dashboard'!`$A`$2:`$A`$6,'Dashboard'!`$B`$2:`$B`$6,1)"
